# Sync automático del tracker (cada 3h)
# Appends 3 new result rows (172-174) to the tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows data: event_id, fecha, jugador_A, jugador_B, pronostico, cuota
$rows = @(
    @{ Row = 172; EventId = "14851659"; Fecha = "2025-10-16"; JugadorA = "Luciano Darderi";  JugadorB = "Shintaro Mochizuki"; Pronostico = "Gana Luciano Darderi"; Cuota = 1.67 },
    @{ Row = 173; EventId = "14851658"; Fecha = "2025-10-16"; JugadorA = "James Duckworth";   JugadorB = "Flavio Cobolli";      Pronostico = "Gana Flavio Cobolli";    Cuota = 1.5  },
    @{ Row = 174; EventId = "14858335"; Fecha = "2025-10-16"; JugadorA = "Kaichi Uchida";      JugadorB = "Jason Jung";          Pronostico = "Gana Jason Jung";        Cuota = 2.63 }
)

foreach ($r in $rows) {
    $n = $r.Row

    # Columns A, B, G and H carry text-typed (inline string) content in this
    # workbook, even when the text looks numeric/date-like (event_id) or is
    # empty (resultado/profit, still pending). Force text format so Excel
    # doesn't silently coerce them into numbers/dates.
    $ws.Range("A$n").NumberFormat = "@"
    $ws.Range("B$n").NumberFormat = "@"
    $ws.Range("G$n").NumberFormat = "@"
    $ws.Range("H$n").NumberFormat = "@"

    $ws.Range("A$n").Value = $r.EventId
    $ws.Range("B$n").Value = $r.Fecha
    $ws.Range("C$n").Value = $r.JugadorA
    $ws.Range("D$n").Value = $r.JugadorB
    $ws.Range("E$n").Value = $r.Pronostico
    $ws.Range("F$n").Value = $r.Cuota

    # resultado / profit: not decided yet for these freshly-synced matches.
    $ws.Range("G$n").Value = ""
    $ws.Range("H$n").Value = ""
}
